# Update the "Total" row of the marksheet: corrected marks obtained (B12)
# and the resulting corr/total display (E12), plus the per-question marking
# value used to recompute the total (B11).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 70
$ws.Range("E12").Value = "70/140"
